$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (shifts "emotions" and everything below it down by one).
# Excel's default row-insert behaviour copies formatting from the row above,
# which matches styles s="3"/"4"/"3" already used by the "location" row.
$ws.Rows.Item(4).Insert()

# Rename "place" -> "location" and broaden its description.
$ws.Range("A3").Value = "location"
$ws.Range("B3").Value = "The place that the user is at the current timestamp, maybe gym, work, car, open-air site, etc"

# Populate the newly inserted row with the "event" feature.
$ws.Range("A4").Value = "event"
$ws.Range("B4").Value = "The particular event that takes place at that moment, may be working, commuting, training, celebration, etc"
$ws.Range("C4").Value = "User-specific"
